$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(163).Insert()

$ws.Range("A163").Value = 5
$ws.Range("B163").Value = "Macroferia Regional de Talca"
$ws.Range("C163").Value = "Maule"
$ws.Range("D163").Value = 45267
$ws.Range("E163").Value = 7
$ws.Range("F163").Value = 100112022
$ws.Range("G163").Value = "Arveja Verde"
$ws.Range("H163").Value = "Sin especificar"
$ws.Range("I163").Value = "Primera"
$ws.Range("J163").Value = 400
$ws.Range("K163").Value = 20000
$ws.Range("L163").Value = 20000
$ws.Range("M163").Value = 20000
$ws.Range("N163").Value = "`$/saco 25 kilos"
$ws.Range("O163").Value = "Región del Maule"
$ws.Range("P163").Value = 800
$ws.Range("Q163").Value = 25
$ws.Range("R163").Value = "Hortaliza"

Write-Host "Done"
